$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-5 with recomputed values (Natmi rerun per Dr Hou advice) ---

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 4.619088000000001
$ws.Range("H2").Value = 13.857264
$ws.Range("K2").Value = 2
$ws.Range("M2").Value = 0.1341725
$ws.Range("N2").Value = 0.268345
$ws.Range("O2").Value = 0.01633512969336317
$ws.Range("P2").Value = 0.01188426112752495
$ws.Range("Q2").Value = 0.6197545846800001
$ws.Range("R2").Value = 3.71852750808
$ws.Range("S2").Value = 0.01633512969336317
$ws.Range("T2").Value = 0.01188426112752495

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 4.619088000000001
$ws.Range("H3").Value = 13.857264
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 5.923689
$ws.Range("N3").Value = 17.771067
$ws.Range("O3").Value = 0.7211927040052828
$ws.Range("P3").Value = 0.787031622511101
$ws.Range("Q3").Value = 27.362040775632
$ws.Range("R3").Value = 246.258366980688
$ws.Range("S3").Value = 0.7211927040052828
$ws.Range("T3").Value = 0.787031622511101

# Row 4
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 4.619088000000001
$ws.Range("H4").Value = 13.857264
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.007277666666666668
$ws.Range("N4").Value = 0.021833
$ws.Range("O4").Value = 0.00088603572911786
$ws.Range("P4").Value = 0.0009669234500260939
$ws.Range("Q4").Value = 0.03361618276800001
$ws.Range("R4").Value = 0.302545644912
$ws.Range("S4").Value = 0.00088603572911786
$ws.Range("T4").Value = 0.0009669234500260939

# Row 5
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 4.619088000000001
$ws.Range("H5").Value = 13.857264
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.039371
$ws.Range("N5").Value = 0.118113
$ws.Range("O5").Value = 0.004793310038624915
$ws.Range("P5").Value = 0.005230899530661476
$ws.Range("Q5").Value = 0.181858113648
$ws.Range("R5").Value = 1.636723022832
$ws.Range("S5").Value = 0.004793310038624915
$ws.Range("T5").Value = 0.005230899530661476

# --- Append two new target-cluster rows: Neutro (row 6) and sCs (row 7) ---

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Wnt5a"
$ws.Range("C6").Value = "Fzd2"
$ws.Range("D6").Value = "Neutro"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 4.619088000000001
$ws.Range("H6").Value = 13.857264
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.1820473333333333
$ws.Range("N6").Value = 0.5461419999999999
$ws.Range("O6").Value = 0.02216375785150397
$ws.Range("P6").Value = 0.02418712530775207
$ws.Range("Q6").Value = 0.840892652832
$ws.Range("R6").Value = 7.568033875487999
$ws.Range("S6").Value = 0.02216375785150397
$ws.Range("T6").Value = 0.02418712530775207

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Wnt5a"
$ws.Range("C7").Value = "Fzd2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 4.619088000000001
$ws.Range("H7").Value = 13.857264
$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 1
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.927182
$ws.Range("N7").Value = 3.854364
$ws.Range("O7").Value = 0.2346290626821072
$ws.Range("P7").Value = 0.1706991680729343
$ws.Range("Q7").Value = 8.901823250016001
$ws.Range("R7").Value = 53.410939500096
$ws.Range("S7").Value = 0.2346290626821072
$ws.Range("T7").Value = 0.1706991680729343
